$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOCU")

$ws.Range("B7").Value = 124000000.0
$ws.Range("C7").Value = 115000000.0
$ws.Range("D7").Value = 64069000.0
$ws.Range("E7").Value = 63157000.0
$ws.Range("F7").Value = 25998000.0
